$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2: Rentabilidade Dia for Constellation - text value changed ("1.96" -> "0,68")
$ws.Range("C2").Value = "'0,68"

# Row 3 (Nucleo): C3/D3 were scraped as numbers, now arrive as text (comma decimal)
$ws.Range("C3").Value = "'0"
$ws.Range("D3").Value = "'-12,02"

# Row 4 (Dynamo): C4/D4 were scraped as numbers, now arrive as text (comma decimal)
$ws.Range("C4").Value = "'-0,01"
$ws.Range("D4").Value = "'-17,71"
